$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 42632.882152777776
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = "Neutral"
$ws.Cells.Item(3, 4).Value = 16
$ws.Cells.Item(3, 5).Value = 8173
$ws.Cells.Item(3, 6).Value = 412
$ws.Cells.Item(3, 7).Value = 60
$ws.Cells.Item(3, 8).Value = 38
$ws.Cells.Item(3, 9).Value = 71
$ws.Cells.Item(3, 10).Value = 28
$ws.Cells.Item(3, 11).Value = 10490
$ws.Cells.Item(3, 12).Value = 80
$ws.Cells.Item(3, 13).Value = 50
$ws.Cells.Item(3, 14).Value = 5
$ws.Cells.Item(3, 15).Value = 2
$ws.Cells.Item(3, 16).Value = "Named"
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 1.77
$ws.Cells.Item(3, 19).Value = 0.1132
$ws.Cells.Item(3, 19).NumberFormat = $ws.Cells.Item(2, 19).NumberFormat
$ws.Cells.Item(3, 20).Value = -4.05
$ws.Cells.Item(3, 21).Value = 5.85
$ws.Cells.Item(3, 22).Value = "N/A"
$ws.Cells.Item(3, 23).Value = 0
